# Auto-generated Excel COM-interop script
# Applies refreshed market-price / profit figures to the Leve profit tables
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Rule observed in the data: column M ("LeveProfitNQ") is present only when
# column K ("LevePriceNQ") is non-zero, and column N ("LeveProfitHQ") is present
# only when column L ("LevePriceHQ") is non-zero. When a refreshed row crosses
# that boundary, the now-inapplicable cell is cleared so it disappears from the
# saved XML (matching the source system's export behaviour), and the newly
# applicable cell is written.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 3142.7144
$ws.Range("J16").Value = 3199.8
$ws.Range("L16").Value = 3199.8
$ws.Range("N16").Value = -3659.8
$ws.Range("H19").Value = 2020.4
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 2450.5
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 2450.5
$ws.Range("M19").Value = -125
$ws.Range("N19").Value = -2800.5
$ws.Range("H50").Value = 200
$ws.Range("J50").Value = 200
$ws.Range("L50").Value = 600
$ws.Range("N50").Value = -1550
$ws.Range("H69").Value = 38507.5
$ws.Range("J69").Value = 38507.5
$ws.Range("L69").Value = 115522.5
$ws.Range("N69").Value = -117270.5
$ws.Range("H72").Value = 38507.5
$ws.Range("J72").Value = 38507.5
$ws.Range("L72").Value = 346567.5
$ws.Range("N72").Value = -355303.5
$ws.Range("H86").Value = 8603.762000000001
$ws.Range("I86").Value = 8855.375
$ws.Range("K86").Value = 8855.375
$ws.Range("M86").Value = -7732.375
$ws.Range("H89").Value = 8603.762000000001
$ws.Range("I89").Value = 8855.375
$ws.Range("K89").Value = 44276.875
$ws.Range("M89").Value = -38660.875
$ws.Range("H100").Value = 5711
$ws.Range("I100").Value = 6416
$ws.Range("K100").Value = 6416
$ws.Range("M100").Value = -5875
$ws.Range("H121").Value = 1999
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1999
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 5997
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -9491
$ws.Range("H125").Value = 3262.7
$ws.Range("I125").Value = 3472.2
$ws.Range("K125").Value = 31249.8
$ws.Range("M125").Value = -28789.8
$ws.Range("H137").Value = 7730.5
$ws.Range("I137").Value = 2448.75
$ws.Range("K137").Value = 7346.25
$ws.Range("M137").Value = -4796.25
$ws.Range("H141").Value = 9386.666999999999
$ws.Range("I141").Value = 9864
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 29592
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -24412
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 11112592
$ws.Range("J32").Value = 250
$ws.Range("K32").Value = 11112592
$ws.Range("L32").Value = 250
$ws.Range("M32").Value = -11112305
$ws.Range("N32").Value = -824
$ws.Range("H37").Value = 366
$ws.Range("I37").Value = 366
$ws.Range("K37").Value = 366
$ws.Range("M37").Value = -93
$ws.Range("H45").Value = 1974.5883
$ws.Range("I45").Value = 1979.1428
$ws.Range("K45").Value = 1979.1428
$ws.Range("M45").Value = -1602.1428
$ws.Range("H74").Value = 9622858
$ws.Range("J74").Value = 22682.25
$ws.Range("L74").Value = 22682.25
$ws.Range("N74").Value = -24430.25
$ws.Range("H77").Value = 9622858
$ws.Range("J77").Value = 22682.25
$ws.Range("L77").Value = 113411.25
$ws.Range("N77").Value = -122147.25
$ws.Range("H102").Value = 8866.061
$ws.Range("J102").Value = 9069.929
$ws.Range("L102").Value = 9069.929
$ws.Range("N102").Value = -12313.929
$ws.Range("H132").Value = 5772.946
$ws.Range("I132").Value = 3312.2856
$ws.Range("K132").Value = 9936.856800000001
$ws.Range("M132").Value = -7406.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2941.9443
$ws.Range("I99").Value = 2383.8333
$ws.Range("K99").Value = 2383.8333
$ws.Range("M99").Value = -885.8332999999998
$ws.Range("H105").Value = 1761.85
$ws.Range("I105").Value = 999.4545000000001
$ws.Range("K105").Value = 999.4545000000001
$ws.Range("M105").Value = 747.5454999999999
$ws.Range("H107").Value = 3333
$ws.Range("I107").Value = 2499.5
$ws.Range("K107").Value = 2499.5
$ws.Range("M107").Value = -579.5
$ws.Range("H134").Value = 72299.53
$ws.Range("I134").Value = 2898.5
$ws.Range("K134").Value = 8695.5
$ws.Range("M134").Value = -6160.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 290.46155
$ws.Range("I22").Value = 293.91666
$ws.Range("K22").Value = 293.91666
$ws.Range("M22").Value = 56.08334000000002
$ws.Range("H134").Value = 288624.9
$ws.Range("I134").Value = 401062.6
$ws.Range("K134").Value = 1203187.8
$ws.Range("M134").Value = -1200652.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1335.6111
$ws.Range("I5").Value = 974.5
$ws.Range("J5").Value = 2599.5
$ws.Range("K5").Value = 2923.5
$ws.Range("L5").Value = 7798.5
$ws.Range("M5").Value = -2811.5
$ws.Range("N5").Value = -8022.5
$ws.Range("H33").Value = 1905.1111
$ws.Range("I33").Value = 2890.4
$ws.Range("J33").Value = 673.5
$ws.Range("K33").Value = 17342.4
$ws.Range("L33").Value = 4041
$ws.Range("M33").Value = -17059.4
$ws.Range("N33").Value = -4607
$ws.Range("H37").Value = 82994.664
$ws.Range("J37").Value = 82994.664
$ws.Range("L37").Value = 248983.992
$ws.Range("N37").Value = -249207.992
$ws.Range("H44").Value = 19598.445
$ws.Range("I44").Value = 19598.445
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 58795.335
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -58397.335
$ws.Range("N44").ClearContents()
$ws.Range("H51").Value = 12737.6
$ws.Range("I51").Value = 8933.875
$ws.Range("J51").Value = 27952.5
$ws.Range("K51").Value = 26801.625
$ws.Range("L51").Value = 83857.5
$ws.Range("M51").Value = -26341.625
$ws.Range("N51").Value = -84777.5
$ws.Range("H98").Value = 595
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H135").Value = 1335.6111
$ws.Range("I135").Value = 974.5
$ws.Range("J135").Value = 2599.5
$ws.Range("K135").Value = 8770.5
$ws.Range("L135").Value = 23395.5
$ws.Range("M135").Value = -6235.5
$ws.Range("N135").Value = -28465.5
$ws.Range("H137").Value = 5168.923
$ws.Range("J137").Value = 3761.3333
$ws.Range("L137").Value = 11283.9999
$ws.Range("N137").Value = -21483.9999
$ws.Range("H140").Value = 233307.39
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4748.25
$ws.Range("I70").Value = 4748.25
$ws.Range("K70").Value = 4748.25
$ws.Range("M70").Value = -4478.25
$ws.Range("H73").Value = 4748.25
$ws.Range("I73").Value = 4748.25
$ws.Range("K73").Value = 4748.25
$ws.Range("M73").Value = -3812.25
$ws.Range("H107").Value = 2131.4443
$ws.Range("I107").Value = 2085.375
$ws.Range("K107").Value = 2085.375
$ws.Range("M107").Value = -165.375
$ws.Range("H113").Value = 3262.8235
$ws.Range("I113").Value = 2822.7144
$ws.Range("K113").Value = 2822.7144
$ws.Range("M113").Value = -652.7143999999998
$ws.Range("H132").Value = 45457816
$ws.Range("I132").Value = 58826940
$ws.Range("J132").Value = 2791.6
$ws.Range("K132").Value = 176480820
$ws.Range("L132").Value = 8374.799999999999
$ws.Range("M132").Value = -176478290
$ws.Range("N132").Value = -13434.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45647.332
$ws.Range("I7").Value = 3264.0527
$ws.Range("K7").Value = 3264.0527
$ws.Range("M7").Value = -3152.0527
$ws.Range("H22").Value = 2879.8
$ws.Range("I22").Value = 3063.3635
$ws.Range("J22").Value = 2375
$ws.Range("K22").Value = 3063.3635
$ws.Range("L22").Value = 2375
$ws.Range("M22").Value = -2768.3635
$ws.Range("N22").Value = -2965
$ws.Range("H27").Value = 2879.8
$ws.Range("I27").Value = 3063.3635
$ws.Range("J27").Value = 2375
$ws.Range("K27").Value = 3063.3635
$ws.Range("L27").Value = 2375
$ws.Range("M27").Value = -2956.3635
$ws.Range("N27").Value = -2589
$ws.Range("H68").Value = 2684.5386
$ws.Range("J68").Value = 2899.8333
$ws.Range("L68").Value = 2899.8333
$ws.Range("N68").Value = -4397.8333
$ws.Range("H71").Value = 2684.5386
$ws.Range("J71").Value = 2899.8333
$ws.Range("L71").Value = 14499.1665
$ws.Range("N71").Value = -21987.1665
$ws.Range("H81").Value = 90000
$ws.Range("J81").Value = 90000
$ws.Range("L81").Value = 90000
$ws.Range("N81").Value = -91996
$ws.Range("H84").Value = 90000
$ws.Range("J84").Value = 90000
$ws.Range("L84").Value = 270000
$ws.Range("N84").Value = -279984
$ws.Range("H126").Value = 45647.332
$ws.Range("I126").Value = 3264.0527
$ws.Range("K126").Value = 9792.158100000001
$ws.Range("M126").Value = -7322.158100000001
$ws.Range("H132").Value = 441707.9
$ws.Range("I132").Value = 436595.47
$ws.Range("J132").Value = 500501
$ws.Range("K132").Value = 1309786.41
$ws.Range("L132").Value = 1501503
$ws.Range("M132").Value = -1307256.41
$ws.Range("N132").Value = -1506563

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 81948
$ws.Range("J93").Value = 77922
$ws.Range("L93").Value = 77922
$ws.Range("N93").Value = -82914
$ws.Range("H132").Value = 1858.8572
$ws.Range("I132").Value = 1901.85
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 5705.549999999999
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -3175.549999999999
$ws.Range("N132").Value = -8057
